$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.464.57'
$ws.Range('E2').Value = '  +0.54%  '

$ws.Range('D3').Value = '2.616.27'
$ws.Range('E3').Value = '  +2.34%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.16%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '589.47'
$ws.Range('E5').Value = '  +3.78%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '149.71'
$ws.Range('E6').Value = '  +2.22%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.998'
$ws.Range('E7').Value = '  -0.19%  '

$ws.Range('E8').Value = '  +2.26%  '

$ws.Range('E9').Value = '  +4.73%  '

$ws.Range('E10').Value = '  +1.86%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.359'
$ws.Range('E12').Value = '  +1.96%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '27.94'
$ws.Range('E13').Value = '  +1.62%  '

$ws.Range('D14').Value = '3.082.24'
$ws.Range('E14').Value = '  +2.33%  '

$ws.Range('D15').Value = '63.439.68'
$ws.Range('E15').Value = '  +0.58%  '

$ws.Range('E16').Value = '  +4.63%  '

$ws.Range('D17').Value = '2.637.87'
$ws.Range('E17').Value = '  +3.22%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '11.53'
$ws.Range('E18').Value = '  +1.01%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '346.86'
$ws.Range('E19').Value = '  +3.36%  '

$ws.Range('E20').Value = '  +3.22%  '

$ws.Range('E21').Value = '  +1.57%  '

$ws.Range('E22').Value = '  -0.11%  '

$ws.Range('E23').Value = '  -3.33%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '66.92'
$ws.Range('E24').Value = '  +2.64%  '

$ws.Range('D25').Value = '2.688.85'
$ws.Range('E25').Value = '  +0.31%  '

$ws.Range('E26').Value = '  +0.54%  '

$ws.Range('E27').Value = '  +0.44%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.26'
$ws.Range('E28').Value = '  +12.92%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.57'
$ws.Range('E29').Value = '  +1.25%  '

$ws.Range('E30').Value = '  +2.12%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.998'
$ws.Range('E31').Value = '  -0.28%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.00'
$ws.Range('E32').Value = '  +8.17%  '

$ws.Range('D33').Value = '0.0₃0841'
$ws.Range('E33').Value = '  +2.73%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '469.91'
$ws.Range('E34').Value = '  +15.38%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.66'
$ws.Range('E35').Value = '  +5.53%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '177.27'
$ws.Range('E36').Value = '  +1.00%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.410'
$ws.Range('E37').Value = '  +2.82%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '19.42'
$ws.Range('E38').Value = '  +2.36%  '

$ws.Range('E39').Value = '  +6.94%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.77'
$ws.Range('E41').Value = '  +1.36%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.00'
$ws.Range('E42').Value = '  +0.07%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '153.93'
$ws.Range('E43').Value = '  +0.63%  '

$ws.Range('E44').Value = '  +2.57%  '

$ws.Range('E45').Value = '  +1.12%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0559'
$ws.Range('E46').Value = '  +6.30%  '

$ws.Range('E47').Value = '  +2.17%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0981'
$ws.Range('E48').Value = '  +2.06%  '

$ws.Range('E49').Value = '  +1.78%  '

$ws.Range('E50').Value = '  +0.55%  '

$ws.Range('E51').Value = '  +0.73%  '
